$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8122.2104
$ws.Range("I116").Value = 10164.923
$ws.Range("J116").Value = 3696.3333
$ws.Range("K116").Value = 10164.923
$ws.Range("L116").Value = 3696.3333
$ws.Range("M116").Value = -6722.923000000001
$ws.Range("N116").Value = -10580.3333

$ws.Range("H136").Value = 39950
$ws.Range("J136").Value = 39950
$ws.Range("L136").Value = 39950
$ws.Range("N136").Value = -50150

$ws.Range("H137").Value = 1128.6316
$ws.Range("I137").Value = 788.72
$ws.Range("J137").Value = 1782.3077
$ws.Range("K137").Value = 2366.16
$ws.Range("L137").Value = 5346.9231
$ws.Range("M137").Value = 183.8400000000001
$ws.Range("N137").Value = -10446.9231

$ws.Range("H139").Value = 70270
$ws.Range("J139").Value = 70270
$ws.Range("L139").Value = 70270
$ws.Range("N139").Value = -80550

$ws.Range("H140").Value = 72995.234
$ws.Range("J140").Value = 94193.336
$ws.Range("L140").Value = 94193.336
$ws.Range("N140").Value = -104553.336

$ws.Range("H141").Value = 3193.5
$ws.Range("I141").Value = 3250.2173
$ws.Range("J141").Value = 3007.1428
$ws.Range("K141").Value = 9750.651899999999
$ws.Range("L141").Value = 9021.428400000001
$ws.Range("M141").Value = -4570.651899999999
$ws.Range("N141").Value = -19381.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4575.86
$ws.Range("I32").Value = 3942.2
$ws.Range("J32").Value = 8166.6
$ws.Range("K32").Value = 3942.2
$ws.Range("L32").Value = 8166.6
$ws.Range("M32").Value = -3655.2
$ws.Range("N32").Value = -8740.6

$ws.Range("H61").Value = 4201.6665
$ws.Range("I61").Value = 4666.9165
$ws.Range("J61").Value = 3271.1667
$ws.Range("K61").Value = 4666.9165
$ws.Range("L61").Value = 3271.1667
$ws.Range("M61").Value = -4454.9165
$ws.Range("N61").Value = -3695.1667

$ws.Range("H74").Value = 1176.2222
$ws.Range("I74").Value = 724.5714
$ws.Range("K74").Value = 724.5714
$ws.Range("M74").Value = 149.4286

$ws.Range("H77").Value = 1176.2222
$ws.Range("I77").Value = 724.5714
$ws.Range("K77").Value = 3622.857
$ws.Range("M77").Value = 745.143

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H136").Value = 4201.6665
$ws.Range("I136").Value = 4666.9165
$ws.Range("J136").Value = 3271.1667
$ws.Range("K136").Value = 14000.7495
$ws.Range("L136").Value = 9813.500100000001
$ws.Range("M136").Value = -11450.7495
$ws.Range("N136").Value = -14913.5001

$ws.Range("H139").Value = 54558
$ws.Range("J139").Value = 54558
$ws.Range("L139").Value = 54558
$ws.Range("N139").Value = -64838

$ws.Range("H141").Value = 61285.715
$ws.Range("J141").Value = 64307.69
$ws.Range("L141").Value = 64307.69
$ws.Range("N141").Value = -74667.69

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 8522.5
$ws.Range("I44").Value = 5045
$ws.Range("J44").Value = 12000
$ws.Range("K44").Value = 5045
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = -4548
$ws.Range("N44").Value = -12994

$ws.Range("H140").Value = 89433.336
$ws.Range("J140").Value = 89433.336
$ws.Range("L140").Value = 89433.336
$ws.Range("N140").Value = -99793.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 842.9
$ws.Range("I16").Value = 836.55554
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 836.55554
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -549.55554
$ws.Range("N16").Value = -1474

$ws.Range("H31").Value = 8788.459999999999
$ws.Range("I31").Value = 2665.6445
$ws.Range("J31").Value = 15348.619
$ws.Range("K31").Value = 2665.6445
$ws.Range("L31").Value = 15348.619
$ws.Range("M31").Value = -2370.6445
$ws.Range("N31").Value = -15938.619

$ws.Range("H34").Value = 8788.459999999999
$ws.Range("I34").Value = 2665.6445
$ws.Range("J34").Value = 15348.619
$ws.Range("K34").Value = 2665.6445
$ws.Range("L34").Value = 15348.619
$ws.Range("M34").Value = -2463.6445
$ws.Range("N34").Value = -15752.619

$ws.Range("H86").Value = 3733.7273
$ws.Range("I86").Value = 3254
$ws.Range("J86").Value = 4761.7144
$ws.Range("K86").Value = 3254
$ws.Range("L86").Value = 4761.7144
$ws.Range("M86").Value = -2131
$ws.Range("N86").Value = -7007.7144

$ws.Range("H89").Value = 3733.7273
$ws.Range("I89").Value = 3254
$ws.Range("J89").Value = 4761.7144
$ws.Range("K89").Value = 16270
$ws.Range("L89").Value = 23808.572
$ws.Range("M89").Value = -10654
$ws.Range("N89").Value = -35040.572

$ws.Range("H113").Value = 842.9
$ws.Range("I113").Value = 836.55554
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 836.55554
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 1333.44446
$ws.Range("N113").Value = -5240

$ws.Range("H134").Value = 777.11365
$ws.Range("I134").Value = 773.63635
$ws.Range("J134").Value = 787.5454999999999
$ws.Range("K134").Value = 2320.90905
$ws.Range("L134").Value = 2362.6365
$ws.Range("M134").Value = 214.0909499999998
$ws.Range("N134").Value = -7432.6365

$ws.Range("H138").Value = 49885.715
$ws.Range("J138").Value = 49885.715
$ws.Range("L138").Value = 49885.715
$ws.Range("N138").Value = -60165.715

$ws.Range("H140").Value = 55300
$ws.Range("J140").Value = 55300
$ws.Range("L140").Value = 55300
$ws.Range("N140").Value = -65660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 583.2069
$ws.Range("I98").Value = 502.33334
$ws.Range("K98").Value = 1507.00002
$ws.Range("M98").Value = -9.00001999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 6000
$ws.Range("J48").Value = 6000
$ws.Range("L48").Value = 6000
$ws.Range("N48").Value = -6970

$ws.Range("H80").Value = 1110005
$ws.Range("I80").Value = 1110005
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1110005
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1109007
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 1110005
$ws.Range("I83").Value = 1110005
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 5550025
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5545033
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 2320.85
$ws.Range("I132").Value = 1856.5588
$ws.Range("J132").Value = 4951.8335
$ws.Range("K132").Value = 5569.6764
$ws.Range("L132").Value = 14855.5005
$ws.Range("M132").Value = -3039.6764
$ws.Range("N132").Value = -19915.5005

$ws.Range("H138").Value = 69033.336
$ws.Range("J138").Value = 69033.336
$ws.Range("L138").Value = 69033.336
$ws.Range("N138").Value = -79313.336

$ws.Range("H140").Value = 89844.5
$ws.Range("J140").Value = 89844.5
$ws.Range("L140").Value = 89844.5
$ws.Range("N140").Value = -100204.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 4800
$ws.Range("I42").Value = 4500
$ws.Range("J42").Value = 5100
$ws.Range("K42").Value = 4500
$ws.Range("L42").Value = 5100
$ws.Range("M42").Value = -3937
$ws.Range("N42").Value = -6226

$ws.Range("H49").Value = 4800
$ws.Range("I49").Value = 4500
$ws.Range("J49").Value = 5100
$ws.Range("K49").Value = 4500
$ws.Range("L49").Value = 5100
$ws.Range("M49").Value = -4353
$ws.Range("N49").Value = -5394

$ws.Range("H136").Value = 3346.3809
$ws.Range("I136").Value = 3316.1428
$ws.Range("J136").Value = 3406.8572
$ws.Range("K136").Value = 9948.428400000001
$ws.Range("L136").Value = 10220.5716
$ws.Range("M136").Value = -7398.428400000001
$ws.Range("N136").Value = -15320.5716

$ws.Range("H138").Value = 58485.7
$ws.Range("J138").Value = 58485.7
$ws.Range("L138").Value = 58485.7
$ws.Range("N138").Value = -68765.7

$ws.Range("H139").Value = 54216.668
$ws.Range("J139").Value = 64060
$ws.Range("L139").Value = 64060
$ws.Range("N139").Value = -74340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5692.5386
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 5000.4287
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 5000.4287
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -6248.4287

$ws.Range("H65").Value = 5692.5386
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 5000.4287
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 25002.1435
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -31242.1435

$ws.Range("H113").Value = 14706446
$ws.Range("J113").Value = 83334216
$ws.Range("L113").Value = 250002648
$ws.Range("N113").Value = -250006988

$ws.Range("H139").Value = 58116.668
$ws.Range("J139").Value = 58116.668
$ws.Range("L139").Value = 58116.668
$ws.Range("N139").Value = -68396.66800000001

$ws.Range("H141").Value = 76051.875
$ws.Range("J141").Value = 76051.875
$ws.Range("L141").Value = 76051.875
$ws.Range("N141").Value = -86411.875

